# Update "想去人数" (want-to-go count) figures that changed in the
# latest scrape (gh-pages output regenerated at commit 456a3b4).
#
# Sheet "展览"   (1st worksheet): F3 242->246, F4 863->867, F6 35->37
# Sheet "全部类型" (4th worksheet): F4 242->246, F5 863->867, F7 35->37

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 246
$wsExhibit.Range("F4").Value = 867
$wsExhibit.Range("F6").Value = 37

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 246
$wsAll.Range("F5").Value = 867
$wsAll.Range("F7").Value = 37
